$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 2: "MA of diff" block (header row 18, index-only rows 19-33) ---
$ws.Cells.Item(18, 2).Value = 5
$ws.Cells.Item(18, 3).Value = 8
$ws.Cells.Item(18, 4).Value = 13
$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(20, 1).Value = 2
$ws.Cells.Item(21, 1).Value = 3
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(25, 1).Value = 7
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(30, 1).Value = 12
$ws.Cells.Item(31, 1).Value = 13
$ws.Cells.Item(32, 1).Value = 14
$ws.Cells.Item(33, 1).Value = 15

# --- Table 3: "historical average" block (header row 35, full data rows 36-50) ---
$ws.Cells.Item(35, 2).Value = 5
$ws.Cells.Item(35, 3).Value = 8
$ws.Cells.Item(35, 4).Value = 13
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = 0.39872538696068099
$ws.Cells.Item(36, 3).Value = 0.63669992960508903
$ws.Cells.Item(36, 4).Value = 0.55338855268251896
$ws.Cells.Item(37, 1).Value = 2
$ws.Cells.Item(37, 2).Value = 0.56789888052188198
$ws.Cells.Item(37, 3).Value = 0.68595166832307097
$ws.Cells.Item(37, 4).Value = 0.69857731869086603
$ws.Cells.Item(38, 1).Value = 3
$ws.Cells.Item(38, 2).Value = 0.58819291266663598
$ws.Cells.Item(38, 3).Value = 0.58554910996656695
$ws.Cells.Item(38, 4).Value = 0.64162764380307902
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = 0.97142857142857097
$ws.Cells.Item(39, 3).Value = 0.78333333333333299
$ws.Cells.Item(39, 4).Value = 0.78333333333333299
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = 0.21709289840870499
$ws.Cells.Item(40, 3).Value = 0.77348633967487401
$ws.Cells.Item(40, 4).Value = 0.65643919332058098
$ws.Cells.Item(41, 1).Value = 6
$ws.Cells.Item(41, 2).Value = 0.80977260855035205
$ws.Cells.Item(41, 3).Value = 0.85890645140646305
$ws.Cells.Item(41, 4).Value = 0.87097017304880398
$ws.Cells.Item(42, 1).Value = 7
$ws.Cells.Item(42, 2).Value = 0.49157186808314002
$ws.Cells.Item(42, 3).Value = 0.77962180531827296
$ws.Cells.Item(42, 4).Value = 0.75003247650684501
$ws.Cells.Item(43, 1).Value = 8
$ws.Cells.Item(43, 2).Value = 0.63842240326949995
$ws.Cells.Item(43, 3).Value = 0.71380437246879802
$ws.Cells.Item(43, 4).Value = 0.71549199073222303
$ws.Cells.Item(44, 1).Value = 9
$ws.Cells.Item(44, 2).Value = 0.48012948391959198
$ws.Cells.Item(44, 3).Value = 0.37657079804439603
$ws.Cells.Item(44, 4).Value = 0.27083863185940699
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = 0.53148189368865995
$ws.Cells.Item(45, 3).Value = 0.93058485984792305
$ws.Cells.Item(45, 4).Value = 0.91700779754205197
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = 0.96803788186524997
$ws.Cells.Item(46, 3).Value = 0.96232866798739702
$ws.Cells.Item(46, 4).Value = 0.96648253525442696
$ws.Cells.Item(47, 1).Value = 12
$ws.Cells.Item(47, 2).Value = 0.35093217328895199
$ws.Cells.Item(47, 3).Value = 0.35682895678221199
$ws.Cells.Item(47, 4).Value = 0.34376888595194
$ws.Cells.Item(48, 1).Value = 13
$ws.Cells.Item(48, 2).Value = 0.69159378750413503
$ws.Cells.Item(48, 3).Value = 0.66284814929316105
$ws.Cells.Item(48, 4).Value = 0.68963495936008801
$ws.Cells.Item(49, 1).Value = 14
$ws.Cells.Item(49, 2).Value = 0.67215121233019404
$ws.Cells.Item(49, 3).Value = 0.73102141062633197
$ws.Cells.Item(49, 4).Value = 0.61726110688411695
$ws.Cells.Item(50, 1).Value = 15
$ws.Cells.Item(50, 2).Value = 0.61775702599888604
$ws.Cells.Item(50, 3).Value = 0.80111173615404996
$ws.Cells.Item(50, 4).Value = 0.82044578058469797

# Restore the selection/active cell left behind by the paste operation
$ws.Range("B19:D33").Select()
